$wb = $excel.ActiveWorkbook

# --- Add "Debts" worksheet after the last existing sheet (Jane) ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDebts = $wb.Worksheets.Add($null, $afterSheet)
$wsDebts.Name = "Debts"

$wsDebts.Range("A1").Value = "name"
$wsDebts.Range("B1").Value = "type"
$wsDebts.Range("C1").Value = "year"
$wsDebts.Range("D1").Value = "term"
$wsDebts.Range("E1").Value = "amount"
$wsDebts.Range("F1").Value = "rate"

$wsDebts.Range("A1:F1").Font.Bold = $true
$wsDebts.Range("E1").NumberFormat = """$""#,##0"

$wsDebts.Range("A1:XFD1").Select() | Out-Null

# --- Add "Fixed Assets" worksheet after "Debts" ---
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFixed = $wb.Worksheets.Add($null, $afterSheet2)
$wsFixed.Name = "Fixed Assets"

$wsFixed.Range("A1").Value = "name"
$wsFixed.Range("B1").Value = "type"
$wsFixed.Range("C1").Value = "basis"
$wsFixed.Range("D1").Value = "value"
$wsFixed.Range("E1").Value = "rate"
$wsFixed.Range("F1").Value = "yod"
$wsFixed.Range("G1").Value = "commission"

$wsFixed.Range("A1:G1").Font.Bold = $true
$wsFixed.Range("C1:D1").NumberFormat = """$""#,##0"

$wsFixed.Range("H10").Select() | Out-Null

# "Fixed Assets" ends up the active/selected tab, matching the target workbook view.
$wsFixed.Activate()
